$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "Max"
$ws.Range("B8").Formula = "=MAX(B2:B6)"

$ws.Range("A9").Value = "Min"
$ws.Range("B9").Formula = "=MIN(B2:B6)"

$ws.Range("A10").Value = "Average"
$ws.Range("B10").Formula = "=AVERAGE(B2:B6)"

$ws.Range("B11").Select()
